# Commit: "add chi 2 to methods list"
#
# The stats_terms sheet is a two-column (A: term variant, B: canonical
# "update" term) lookup list kept sorted by column A. A new row is being
# added for the term "chi 2", mapping to the same canonical value as the
# other chi-square variants ("chi-square"). Because the sheet is sorted,
# "chi 2" belongs immediately above the existing "chi square" row (row 25),
# pushing that row and everything below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("stats_terms")

# Insert a brand-new row at 25 (old row 25 "chi square" / "chi-square" and
# everything after it shifts down to make room) and fill in the new term.
$ws.Rows.Item(25).Insert()
$ws.Cells.Item(25, 1).Value = "chi 2"
$ws.Cells.Item(25, 2).Value = "chi-square"

# Keep the sheet's recorded sort range/state in sync with the extra row
# (was A2:B146, now covers one more row -> A2:B147).
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($ws.Range("A2:B147"))
$ws.Sort.Header = 0
$ws.Sort.Apply()

# Leave the selection where the editor ended up after typing the new row.
$ws.Cells.Item(37, 2).Select()
